$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.586.85'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '2.312.28'
$ws.Range("E3").Value = '  +4.01%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '268.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0943'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.10%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").Value = '2.664.09'
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.52'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.14%  '
$ws.Range("E16").Value = '  +8.85%  '
$ws.Range("D17").Value = '2.324.74'
$ws.Range("E17").Value = '  +4.68%  '
$ws.Range("D18").Value = '43.523.99'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000110'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.45%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.52'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.11%  '
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '172.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0894'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.22%  '
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0356'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.91%  '
$ws.Range("E37").Value = '  -2.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("E40").Value = '  +8.67%  '
$ws.Range("E41").Value = '  +11.72%  '
$ws.Range("E42").Value = '  +17.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.102'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '99.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("E49").Value = '  +2.98%  '
$ws.Range("D50").Value = '2.544.82'
$ws.Range("E50").Value = '  +3.99%  '
$ws.Range("E51").Value = '  +13.84%  '
